# Generate Report for handback
#
# Populates the "Latest Target File" / "Latest Handback File" columns (E, F)
# for the two entry rows on both locale sheets, stamps the "Latest Handback
# DateTime" column (G) with the real handback timestamp (replacing the
# zero-date placeholder), and updates the Status column (B) to reflect that
# the files have been handed back and are back in sync with en-US.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column (B2:B3) -> handed back
$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

# Row 2 (da6619c9-...md) - Latest Target File / Latest Handback File
$zh.Hyperlinks.Add(
    $zh.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/66ac320473827d61c81aef5051d3d4976f3e399d/e2e/da6619c9-990a-4ecd-902a-a168e071c4e4.md",
    "",
    "",
    "da6619c9-990a-4ecd-902a-a168e071c4e4.md"
) | Out-Null
$zh.Range("E2").Style = "Hyperlink"

$zh.Hyperlinks.Add(
    $zh.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e6b872c9d45c69cecce6d31bcbdadb47c45f3896/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/da6619c9-990a-4ecd-902a-a168e071c4e4.025fbc935625f4e3cea20fb16dace02e87a470e0.zh-cn.xlf",
    "",
    "",
    "da6619c9-990a-4ecd-902a-a168e071c4e4.025fbc935625f4e3cea20fb16dace02e87a470e0.zh-cn.xlf"
) | Out-Null
$zh.Range("F2").Style = "Hyperlink"

# Latest Handback DateTime for row 2
$zh.Range("G2").Value = "2016-02-17 03:29:52"

# Row 3 (f22a0950-...md) - Latest Target File / Latest Handback File
$zh.Hyperlinks.Add(
    $zh.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/66ac320473827d61c81aef5051d3d4976f3e399d/e2e/f22a0950-cc5c-4031-9b15-468b14d88c9b.md",
    "",
    "",
    "f22a0950-cc5c-4031-9b15-468b14d88c9b.md"
) | Out-Null
$zh.Range("E3").Style = "Hyperlink"

$zh.Hyperlinks.Add(
    $zh.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e6b872c9d45c69cecce6d31bcbdadb47c45f3896/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f22a0950-cc5c-4031-9b15-468b14d88c9b.2787c0ac0eaeb4e20cedc016e532fb2fa5eaf909.zh-cn.xlf",
    "",
    "",
    "f22a0950-cc5c-4031-9b15-468b14d88c9b.2787c0ac0eaeb4e20cedc016e532fb2fa5eaf909.zh-cn.xlf"
) | Out-Null
$zh.Range("F3").Style = "Hyperlink"

# Latest Handback DateTime for row 3
$zh.Range("G3").Value = "2016-02-17 03:29:52"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Status column (B2:B3) -> handed back
$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

# Row 2 (da6619c9-...md) - Latest Target File / Latest Handback File
$de.Hyperlinks.Add(
    $de.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/66ac320473827d61c81aef5051d3d4976f3e399d/e2e/da6619c9-990a-4ecd-902a-a168e071c4e4.md",
    "",
    "",
    "da6619c9-990a-4ecd-902a-a168e071c4e4.md"
) | Out-Null
$de.Range("E2").Style = "Hyperlink"

$de.Hyperlinks.Add(
    $de.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/123832bbfb8cbab1e9cb24b8ab4b60d928e51a08/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/da6619c9-990a-4ecd-902a-a168e071c4e4.025fbc935625f4e3cea20fb16dace02e87a470e0.de-de.xlf",
    "",
    "",
    "da6619c9-990a-4ecd-902a-a168e071c4e4.025fbc935625f4e3cea20fb16dace02e87a470e0.de-de.xlf"
) | Out-Null
$de.Range("F2").Style = "Hyperlink"

# Latest Handback DateTime for row 2
$de.Range("G2").Value = "2016-02-17 03:30:11"

# Row 3 (f22a0950-...md) - Latest Target File / Latest Handback File
$de.Hyperlinks.Add(
    $de.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/66ac320473827d61c81aef5051d3d4976f3e399d/e2e/f22a0950-cc5c-4031-9b15-468b14d88c9b.md",
    "",
    "",
    "f22a0950-cc5c-4031-9b15-468b14d88c9b.md"
) | Out-Null
$de.Range("E3").Style = "Hyperlink"

$de.Hyperlinks.Add(
    $de.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/123832bbfb8cbab1e9cb24b8ab4b60d928e51a08/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f22a0950-cc5c-4031-9b15-468b14d88c9b.2787c0ac0eaeb4e20cedc016e532fb2fa5eaf909.de-de.xlf",
    "",
    "",
    "f22a0950-cc5c-4031-9b15-468b14d88c9b.2787c0ac0eaeb4e20cedc016e532fb2fa5eaf909.de-de.xlf"
) | Out-Null
$de.Range("F3").Style = "Hyperlink"

# Latest Handback DateTime for row 3
$de.Range("G3").Value = "2016-02-17 03:30:11"

"Report generated for handback"
